# Rewrites the prologue chapter per the commit: renames the opening
# Heading1/bookmark, inserts the front-matter sections (by-line,
# Backpage, Dedication, Quote) ahead of it, replaces the old 'Ages'
# copy with the edited 'Introduction' copy, and appends the new
# 'Looking Up to Political Philosophy' / 'Your Politics' / 'Read Now
# or Read Later' sections plus the three-books blurb.
#
# Implemented as a single Range.InsertXML() of the full target body
# (valid Word COM: Range.InsertXML accepts a WordprocessingML
# pkg:package payload) rather than dozens of piecemeal Find/Replace
# calls, since nearly every paragraph in the chapter changed.
$d = $word.ActiveDocument
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:bookmarkStart w:id="20" w:name="techvolution-a-new-philosophy"/><w:r><w:t xml:space="preserve">TECHVOLUTION: A NEW PHILOSOPHY</w:t></w:r><w:bookmarkEnd w:id="20"/></w:p><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">by Johnny Majic</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">December 2019</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:bookmarkStart w:id="21" w:name="backpage"/><w:r><w:t xml:space="preserve">Backpage</w:t></w:r><w:bookmarkEnd w:id="21"/></w:p><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">Are you happy? On the outside Westerners have flashy products, full employment, and stable politics. On the inside our middle-class deals with stagnate wages, crushing debts, increasing inequality, antisocial behavior, political hatred, rat race lifestyles, climate change. We have the freedom to vote any ideology, and we’ve tried them all, yet our problems keep increasing. In this book I argue the middle-class needs a totally new philosophy, one that says we either step up to solve our own problems, or keep suffering the unknowns of a degenerating society.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:bookmarkStart w:id="22" w:name="dedication"/><w:r><w:t xml:space="preserve">Dedication</w:t></w:r><w:bookmarkEnd w:id="22"/></w:p><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">This book is for the people in the open-source community. For the unspoken philosophy you handed down to me. For selflessly solving problems when updating Wikipedia, writing free software, and building cool new tech for us all. For showing me who the real protagonists are. You gave the best education a political guy can get. I hope I’ve learned enough to write this book.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">For the antagonists, the self-entitled delegators keeping Blockbuster—and all it stands for—alive, thanks for causing history to repeat itself yet again. I hope you learn something.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">This is especially dedicated to the extras. Because when it’s time for you to act, I pray you pick the right side.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:bookmarkStart w:id="23" w:name="quote"/><w:r><w:t xml:space="preserve">Quote</w:t></w:r><w:bookmarkEnd w:id="23"/></w:p><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">We are really the revolutionaries in the world today—not the kids with long hair and beards who were wrecking the schools a few years ago.</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Gordon Moore, co-founder of Intel, 1973.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:bookmarkStart w:id="24" w:name="introduction"/><w:r><w:t xml:space="preserve">Introduction</w:t></w:r><w:bookmarkEnd w:id="24"/></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:bookmarkStart w:id="25" w:name="looking-up"/><w:r><w:t xml:space="preserve">Looking Up</w:t></w:r><w:bookmarkEnd w:id="25"/></w:p><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">How do we develop through technological Ages? Going from the Stone Age to the Bronze Age or the Medieval Age to Scientific Age happens over generations—it’s hard to get our head around. I’ve been thinking about it for years and only now have a grasp. There’s a lot to think about when a civilization’s entire way of life upgrades to new technology.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Play games like Age of Empires, Anno 1800, and Civilization, and you’ll get a good idea of how civilizations develop. You’ll be the player</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">god</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">that saves up the resources, invests in the science, and decides when to click</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">upgrade.</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">When you do, civilization instantly transforms. Like magic, every house, factory, building, road, library, school is suddenly new and improved.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">The above is just a game, of course. But we’re all in there somewhere, trying to figure out how to earn a living while our society transforms around us. It’s hard—there is no player god and no magic button. Instead, we, the people inside the game, must collect the resources, do the research, and manually upgrade every part of our civilization bit by bit with blood, sweat, and tears. We can’t sit back, click, and watch. We must live it.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Improving technology is hard work, but at least we understand how to do research and development. Improving politics isn’t so straightforward. Politics is contentious, especially when society has one foot in an old and new era. This split identity is the world of politics and philosophy, civil wars, and revolutions I want to help you make a better sense of. Because in 1970, we entered the Digital Age. Since then, digital innovators, geeks, and hackers have been upgrading our society with digital tools—while the powers that be have been holding them back.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">And that’s the reason life for the middle-class has gotten worse. Things will only improve once we finish the Digital Age upgrade. There’s no game, cheats, or gimmicks. There’s just you and your friends on Main Street. You all have to get out the sweat-bands, work gloves, and new ideas to upgrade every piece of our civilization, bit by bit.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">The job before us is as natural as evolution. Life only adapts and diversifies when new and improved lifeforms replace old ones. Likewise, right now, our civilization is moving from the Industrial Age to the Digital Age. You can’t sit out. Either you help our evolution, or old it back.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">This book integrates technology, history, politics, and philosophy, but it’s not written for any buff. It uses the power of stories to help the people inside today’s game—the middle-class—improve their own lives. The first step is a new mindset. Because we already have digital tools. We only need a modern philosophy to enjoy all the wonders of the Digital Age.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:bookmarkStart w:id="26" w:name="looking-up-to-political-philosophy"/><w:r><w:t xml:space="preserve">Looking Up to Political Philosophy</w:t></w:r><w:bookmarkEnd w:id="26"/></w:p><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">This book’s illustration has a regular person standing on Main Street. His or her’s job is to upgrade society to a new Age. What does that mean inside one person’s life? So she’s looking up, hoping for guidance and sees two options—the heroes and the villains—fighting each other. In the sky high above is the</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">player god</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">who determines the winner.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">In the computer game, when the player god presses the upgrade button, every character gets all-new tools: quills become pens, telephones become smartphones, and so on. It’s what we do when shopping, reading reviews, and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">saving up</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">for a new cool thing. Our shopping determines the health of our entire civilization. Yet we don’t think of ourselves as being</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">player gods</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r><w:r><w:t xml:space="preserve">; but in fact, it’s all our individual choices that make up the</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">player god.</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">The more individuals realize the</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">big picture</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and act like their actions matter, the healthier their society. And the job of giving them perspective is the world of political philosophy. And it’s a mysterious world with one clear goal; make people act with the big picture in mind. The regular person looking up from Main Street is interpreting his society’s ruling philosophy, and with her interpretation, she decides how he should act within the game.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">But what the heck is the big picture, and who knows it best? In everyday life, everyone says their fighting for good—no one says they’re evil. But looking back at history and the villains and heroes easily pop out. Both must exist today. But who is who?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">So there you are, on Main Street, looking up, trying to figure if Donald Trump is evil, or if the United Nations is trying to make a world government. You ask if Fox News is manufacturing fake news, or maybe</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">all of them</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">are. If so, what’s real news? Should the government-run health insurance, or a big corporation?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">It’s all very confusing and emotional. The sad part is you have to pick sides. After all, today, you’re the one supposed to upgrade civilization to the Digital Age. Wouldn’t it be cool if your political philosophy showed you the difference between good and evil? It is cool; it’s calming to understand politics. And that’s one reason why I wrote this book for you.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:bookmarkStart w:id="27" w:name="your-politics"/><w:r><w:t xml:space="preserve">Your Politics</w:t></w:r><w:bookmarkEnd w:id="27"/></w:p><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">When a computer player</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">upgrades</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">to a new Age, we can see technology changed. Damning a river, or putting up telephone poles is very obvious.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">When the computer player</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">upgrades</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">the game’s philosophy, he installs a new ruling philosophy. For example, going from Divine Monarchy to American democracy. But we can’t see this change in real-world objects. The change happens in our minds. And that’s much harder to see. When the upgrade happens, every person on Main Street looks up and asks,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">What does American democracy mean? Am I supposed to be conservative or progressive? Libertarian or socialist? Republican or Democrat?</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">These Factions are interpretations of the ruling philosophy—in this case, American democracy. A person uses the factions to decide what they think the</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">big picture</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and thus form their own political opinions.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">People are good at expressing and fighting for their political opinions. But understanding where their mindset comes from is entirely different. Nobody ever says,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">My family is old money, so I’m coasting through life. I only work when I have to protect my entitlement; what most people would call undeserved inheritance.</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Instead, the person says,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">White privilege doesn’t exist stupid!</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Nobody says,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">I’m a gullible person who works with lots of Urban Studies graduates, so I say xher, and I don’t like Donald Trump because he’s mean, vulgar, and has orange hair!</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">instead that person says;</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">Trump is not presidential.</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">The point is, the real reason for your political opinions is hard to know. We all look up at the same ruling philosophy; we all try to</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">think of the big picture,</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">yet we arrive at different answers. Our opinions are a byproduct of our psychology, job, skillset, family, culture, ethnicity, religion, and time-period. Who can understand how all these factors intermix to make your opinion? Why does one parent cry tears of sadness, the other of joy, on the same election night?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">I can’t explain it either. But I can explain where you are in the story of civilization. Because starting a new Age always has the same cast of characters. Some people fight to retain power for themselves — people who want to distribute power to the people. And regular people on Main Street who either help the good guys or don’t. Whatever the reason for not getting involved doesn’t matter. People either help the good guys upgrade civilization, or don’t.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Today this translates to the digital technologists fighting to distribute the power of new digital technology to the people on Main Street. The people on Main Street have a choice, do they step up to help, or not.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:bookmarkStart w:id="28" w:name="read-now-or-read-later"/><w:r><w:t xml:space="preserve">Read Now or Read Later</w:t></w:r><w:bookmarkEnd w:id="28"/></w:p><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">There’s a heck of a lot of emotions in Western politics today. Better put—people fucking hate each other. And few of those people know why. They’re just mad. And not mad in the</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">I stepped in the puddle</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">way. Mad in the comic book supervillain way.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Let’s get some perspective. Democracy spread around the world on the back of the Industrial Age. Conservatism and Progressivism were the biggest factions of that ruling philosophy. Now, at the start of the Digital Age, rather than rethink hundreds of years of</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">left vs. right</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">politics, each side blames each other for their shared outdatedness. We can make a comparison to past feuds within old ruling philosophies (the Thirty Years War comes to my mind). Either way, impeachments, disputed elections, rallies, and political deadlock is a fucking horrible place for our society to be in. When we should be upgrading our ruling philosophy to the Digital Age, the rantings and ravings of outdated ideologues disengage ever-more Westerners from politics.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">We achieved great things in the Industrial Age. It wasn’t perfect, but new tools in medicine, physics, sanitation, education, agriculture empowered us to fight disease, ignorance, and hunger. But we’re not in the Industrial Age anymore.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">There’s a good reason why Civilization computer games are about players building new cities, farms, sewers, and roads. These are the things that make a civilization great. Political philosophy is the way to let the people see the big picture and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">play the game</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">well; it’s not trolling your neighbor. Today the big picture is we need a new philosophy. One that isn’t spoken by battle-worn ideologues, or doesn’t scare regular people away with</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">the other side is evil.</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">We need one that looks ahead towards the boundless horizon of a new Digital Age. One that sees Tesla, SpaceX, Blue Origin, as a sneak peek for the wondrous Digital Age future we’re all supposed to build bit by bit.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">I could have written this has a political theory book or a history book. But I didn’t (I felt forced to write it actually). It’s written to the people on Main Street right now. Particularly the ones sick of politics but who love checking out what Elon is up. People who’ve played Anno 1800, SimCity, or Civilization and have been the player god, and who realize it’s the actions of regular people that make up the player god. People who can understand most of today’s political problems, like road traffic, hospital wait-times, expensive education, are engineering problems, not political problems.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">We’ll get there quickly and without too much bloodshed. But that’s only if the people on Main Street realize the big picture, that we’re all the player god. And if we don’t step up, at least we’ll be able to read what we could have done when regretting the past in a big picture of remorse.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">===</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Techvolution has three books:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">A New Philosophy</w:t></w:r><w:r><w:t xml:space="preserve">: Upgrades our political philosophy to the Digital Age. This book empowers digital protagonists like Gates, Torvalds, and Musk with the mass support of Main Street.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Building the Life Star</w:t></w:r><w:r><w:t xml:space="preserve">: Explains the new political and economic truth we’ll see in the Digital Age. This book guides us to more accurate definitions of right and wrong that’ll help us remove most political uncertainty in society. Hint the Life Star will do for politics and economics what the scientific instruments did for physics and</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Starting Starfleet</w:t></w:r><w:r><w:t xml:space="preserve">: Shows the wondrous world of perpetual profit awaiting us when we stop wasting our money, keeping the Industrial Age alive. In this part, we’ll outline the digital products that’ll make the Digital Age way of life, where collaboration replaces competition as the default political and economical answer. We’ll send the massive saving to Musk and his protegees and hence start Starfleet.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Content.InsertXML($xml)
